$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column at D (old D->E, old E->F)
$ws.Columns.Item(4).Insert()

# 2. Adjust column widths (C widened, new D narrower)
$ws.Columns.Item(3).ColumnWidth = 24.666666666666668
$ws.Columns.Item(4).ColumnWidth = 13.0

# 3. Fix up hyperlinks: they do not auto-shift with the column insert, so
#    rebuild them pointing at the new F column, preserving order/targets.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F10"), "http://www.boendgen.com/tipps-tricks.html")
$ws.Hyperlinks.Add($ws.Range("F9"), "http://www.bauzentrum-widmann.de/")
$ws.Hyperlinks.Add($ws.Range("F11"), "http://www.bauexpertenforum.de/")
$ws.Hyperlinks.Add($ws.Range("F12"), "http://www.bau.net/forum/tiefbau/11045.php")

# 4. New header for inserted column
$ws.Range("D1").Value = "Kosten"

# 5. Update the "Zement" row (row 14) content across C/D/E
$ws.Range("C14").Value = "ca.  2,4t = 96 Sack"
$ws.Range("D14").Value = "96 * ~3 = 288€"
$ws.Range("D14").Characters(6,1).Font.Size = 10
$ws.Range("D14").Characters(7,8).Font.Size = 11
$ws.Range("E14").Value = "4 m³ x 2 = 8 m³"

# 6. Update the "Kies(Beton)" row (row 15): swap the two remaining values
$ws.Range("C15").Value = "ca 2t/m³ "
$ws.Range("F15").Value = "16t"

# 7. Copy formatting into the two new empty cells that need styling
$ws.Range("C44").Copy()
$ws.Range("D44").PasteSpecial(-4122)
$ws.Range("C86").Copy()
$ws.Range("D86").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 8. Move the selection to match the saved view
$ws.Range("D16").Select()
